$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.794650418452477
$ws.Cells.Item(2, 3).Value = 0.1649669270842082
$ws.Cells.Item(2, 4).Value = 0.07906294041258377
$ws.Cells.Item(2, 5).Value = 0.1092766512264234
$ws.Cells.Item(2, 7).Value = 0.002408672199012024
$ws.Cells.Item(2, 13).Value = 0.3398181348944291
$ws.Cells.Item(2, 14).Value = 1.078372830862985
$ws.Cells.Item(2, 15).Value = 2.571205963064813
$ws.Cells.Item(3, 2).Value = 0.7050580247677942
$ws.Cells.Item(3, 3).Value = 0.1443711910447973
$ws.Cells.Item(3, 4).Value = 0.07163472829491013
$ws.Cells.Item(3, 5).Value = 0.103144775336709
$ws.Cells.Item(3, 7).Value = 0.002412187661532699
$ws.Cells.Item(3, 13).Value = 0.3047588530634329
$ws.Cells.Item(3, 14).Value = 1.094194201261868
$ws.Cells.Item(3, 15).Value = 2.536162412778424
$ws.Cells.Item(4, 2).Value = 0.6501349642453818
$ws.Cells.Item(4, 3).Value = 0.131673953113534
$ws.Cells.Item(4, 4).Value = 0.06711059129069952
$ws.Cells.Item(4, 5).Value = 0.09947083137834767
$ws.Cells.Item(4, 7).Value = 0.002414460969516221
$ws.Cells.Item(4, 13).Value = 0.2833408364309022
$ws.Cells.Item(4, 14).Value = 1.104412393761754
$ws.Cells.Item(4, 15).Value = 2.516650049453688
$ws.Cells.Item(5, 2).Value = 0.6277756577867706
$ws.Cells.Item(5, 3).Value = 0.1264868883372401
$ws.Cells.Item(5, 4).Value = 0.06527616491021604
$ws.Cells.Item(5, 5).Value = 0.09799630466626752
$ws.Cells.Item(5, 7).Value = 0.002415416320107952
$ws.Cells.Item(5, 13).Value = 0.2746398735687592
$ws.Cells.Item(5, 14).Value = 1.108702923746931
$ws.Cells.Item(5, 15).Value = 2.509200789869368
$ws.Cells.Item(6, 2).Value = 0.6240642701172305
$ws.Cells.Item(6, 3).Value = 0.125624807282037
$ws.Cells.Item(6, 4).Value = 0.06497211345470078
$ws.Cells.Item(6, 5).Value = 0.09775282014663844
$ws.Cells.Item(6, 7).Value = 0.002415576707160749
$ws.Cells.Item(6, 13).Value = 0.2731967124328136
$ws.Cells.Item(6, 14).Value = 1.109423000906313
$ws.Cells.Item(6, 15).Value = 2.507994122450157
$ws.Cells.Item(7, 2).Value = 0.6498333280427744
$ws.Cells.Item(7, 3).Value = 0.1316040504109139
$ws.Cells.Item(7, 4).Value = 0.06708581441951367
$ws.Cells.Item(7, 5).Value = 0.09945085408926957
$ws.Cells.Item(7, 7).Value = 0.002414473736159958
$ws.Cells.Item(7, 13).Value = 0.2832233829752653
$ws.Cells.Item(7, 14).Value = 1.104469745225749
$ws.Cells.Item(7, 15).Value = 2.516547555279175
$ws.Cells.Item(8, 2).Value = 0.7637411090241244
$ws.Cells.Item(8, 3).Value = 0.1578762217281167
$ws.Cells.Item(8, 4).Value = 0.07649401258021271
$ws.Cells.Item(8, 5).Value = 0.1071433232868699
$ws.Cells.Item(8, 7).Value = 0.002409860559235481
$ws.Cells.Item(8, 13).Value = 0.3277070162919884
$ws.Cells.Item(8, 14).Value = 1.083723383073734
$ws.Cells.Item(8, 15).Value = 2.558705766227575
$ws.Cells.Item(9, 2).Value = 0.9877985043161175
$ws.Cells.Item(9, 3).Value = 0.208987213495675
$ws.Cells.Item(9, 4).Value = 0.09523932166669624
$ws.Cells.Item(9, 5).Value = 0.1229633159171897
$ws.Cells.Item(9, 7).Value = 0.002401720814845593
$ws.Cells.Item(9, 13).Value = 0.4158156889612314
$ws.Cells.Item(9, 14).Value = 1.047043630012491
$ws.Cells.Item(9, 15).Value = 2.657375315305217
$ws.Cells.Item(10, 2).Value = 1.152840078624138
$ws.Cells.Item(10, 3).Value = 0.2462916621214504
$ws.Cells.Item(10, 4).Value = 0.1091984699683763
$ws.Cells.Item(10, 5).Value = 0.1350531121097021
$ws.Cells.Item(10, 7).Value = 0.002396287393241562
$ws.Cells.Item(10, 13).Value = 0.481111306763168
$ws.Cells.Item(10, 14).Value = 1.022543753807662
$ws.Cells.Item(10, 15).Value = 2.739760306137043
$ws.Cells.Item(11, 2).Value = 1.228017552440633
$ws.Cells.Item(11, 3).Value = 0.263209583560581
$ws.Cells.Item(11, 4).Value = 0.1155909531735375
$ws.Cells.Item(11, 5).Value = 0.140658571091997
$ws.Cells.Item(11, 7).Value = 0.002393933077300165
$ws.Cells.Item(11, 13).Value = 0.510944495146191
$ws.Cells.Item(11, 14).Value = 1.011931606832801
$ws.Cells.Item(11, 15).Value = 2.779418439190124
$ws.Cells.Item(12, 2).Value = 1.256499455620258
$ws.Cells.Item(12, 3).Value = 0.2696084537019772
$ws.Cells.Item(12, 4).Value = 0.1180178069738815
$ws.Cells.Item(12, 5).Value = 0.1427967296401178
$ws.Cells.Item(12, 7).Value = 0.002393058342100301
$ws.Cells.Item(12, 13).Value = 0.5222606037722102
$ws.Cells.Item(12, 14).Value = 1.007989890091572
$ws.Cells.Item(12, 15).Value = 2.79475182300348
$ws.Cells.Item(13, 2).Value = 1.25036476026969
$ws.Cells.Item(13, 3).Value = 0.2682306810977195
$ws.Cells.Item(13, 4).Value = 0.117494865658017
$ws.Cells.Item(13, 5).Value = 0.1423355453124842
$ws.Cells.Item(13, 7).Value = 0.002393245986296093
$ws.Cells.Item(13, 13).Value = 0.5198226314678465
$ws.Cells.Item(13, 14).Value = 1.008835386606327
$ws.Cells.Item(13, 15).Value = 2.791435429388798
$ws.Cells.Item(14, 2).Value = 1.230360501194298
$ws.Cells.Item(14, 3).Value = 0.2637361748916192
$ws.Cells.Item(14, 4).Value = 0.1157904880988525
$ws.Cells.Item(14, 5).Value = 0.1408341665358108
$ws.Cells.Item(14, 7).Value = 0.002393860776299139
$ws.Cells.Item(14, 13).Value = 0.5118750974266959
$ws.Cells.Item(14, 14).Value = 1.011605777422689
$ws.Cells.Item(14, 15).Value = 2.780673586464502
$ws.Cells.Item(15, 2).Value = 1.218109102314827
$ws.Cells.Item(15, 3).Value = 0.2609821702436363
$ws.Cells.Item(15, 4).Value = 0.1147473115936748
$ws.Cells.Item(15, 5).Value = 0.1399165548414771
$ws.Cells.Item(15, 7).Value = 0.002394239537170616
$ws.Cells.Item(15, 13).Value = 0.5070094757834482
$ws.Cells.Item(15, 14).Value = 1.013312739923022
$ws.Cells.Item(15, 15).Value = 2.774122824578285
$ws.Cells.Item(16, 2).Value = 1.147929012600457
$ws.Cells.Item(16, 3).Value = 0.2451849805703148
$ws.Cells.Item(16, 4).Value = 0.1087815646286714
$ws.Cells.Item(16, 5).Value = 0.1346889345608204
$ws.Cells.Item(16, 7).Value = 0.002396443607893556
$ws.Cells.Item(16, 13).Value = 0.4791642726683563
$ws.Cells.Item(16, 14).Value = 1.023248027040239
$ws.Cells.Item(16, 15).Value = 2.737212617444072
$ws.Cells.Item(17, 2).Value = 1.104900898416645
$ws.Cells.Item(17, 3).Value = 0.2354805156462305
$ws.Cells.Item(17, 4).Value = 0.105132675336975
$ws.Cells.Item(17, 5).Value = 0.1315092469543799
$ws.Cells.Item(17, 7).Value = 0.002397825735705983
$ws.Cells.Item(17, 13).Value = 0.4621155859539954
$ws.Cells.Item(17, 14).Value = 1.029479673963781
$ws.Cells.Item(17, 15).Value = 2.715129454414637
$ws.Cells.Item(18, 2).Value = 1.080161631652913
$ws.Cells.Item(18, 3).Value = 0.2298938695033996
$ws.Cells.Item(18, 4).Value = 0.1030379175502389
$ws.Cells.Item(18, 5).Value = 0.1296903112169829
$ws.Cells.Item(18, 7).Value = 0.002398631752038974
$ws.Cells.Item(18, 13).Value = 0.4523218424132693
$ws.Cells.Item(18, 14).Value = 1.033114103795361
$ws.Cells.Item(18, 15).Value = 2.702632942032238
$ws.Cells.Item(19, 2).Value = 1.071786963362683
$ws.Cells.Item(19, 3).Value = 0.2280014893366911
$ws.Cells.Item(19, 4).Value = 0.1023293513930241
$ws.Cells.Item(19, 5).Value = 0.129076148740161
$ws.Cells.Item(19, 7).Value = 0.002398906556375837
$ws.Cells.Item(19, 13).Value = 0.4490079386599888
$ws.Cells.Item(19, 14).Value = 1.034353265839822
$ws.Cells.Item(19, 15).Value = 2.698437006214249
$ws.Cells.Item(20, 2).Value = 1.109480348476779
$ws.Cells.Item(20, 3).Value = 0.2365140797604681
$ws.Cells.Item(20, 4).Value = 0.105520692758418
$ws.Cells.Item(20, 5).Value = 0.1318466996754353
$ws.Cells.Item(20, 7).Value = 0.002397677462679362
$ws.Cells.Item(20, 13).Value = 0.4639291822577718
$ws.Cells.Item(20, 14).Value = 1.028811112360403
$ws.Cells.Item(20, 15).Value = 2.717458998410109
$ws.Cells.Item(21, 2).Value = 1.236235865281117
$ws.Cells.Item(21, 3).Value = 0.2650565267560978
$ws.Cells.Item(21, 4).Value = 0.1162909374455126
$ws.Cells.Item(21, 5).Value = 0.1412747351137114
$ws.Cells.Item(21, 7).Value = 0.002393679742565769
$ws.Cells.Item(21, 13).Value = 0.5142089651200195
$ws.Cells.Item(21, 14).Value = 1.010789958291811
$ws.Cells.Item(21, 15).Value = 2.783826015427735
$ws.Cells.Item(22, 2).Value = 1.319158307144505
$ws.Cells.Item(22, 3).Value = 0.2836664444140808
$ws.Cells.Item(22, 4).Value = 0.1233658443538985
$ws.Cells.Item(22, 5).Value = 0.1475269399043242
$ws.Cells.Item(22, 7).Value = 0.002391164846005157
$ws.Cells.Item(22, 13).Value = 0.5471801535919809
$ws.Cells.Item(22, 14).Value = 0.9994602788932507
$ws.Cells.Item(22, 15).Value = 2.829041883420473
$ws.Cells.Item(23, 2).Value = 1.274893829047926
$ws.Cells.Item(23, 3).Value = 0.2737380577900979
$ws.Cells.Item(23, 4).Value = 0.1195865251153378
$ws.Cells.Item(23, 5).Value = 0.1441816503323281
$ws.Cells.Item(23, 7).Value = 0.00239249816886279
$ws.Cells.Item(23, 13).Value = 0.5295726204190885
$ws.Cells.Item(23, 14).Value = 1.005466054637171
$ws.Cells.Item(23, 15).Value = 2.804740153528314
$ws.Cells.Item(24, 2).Value = 1.107409984778599
$ws.Cells.Item(24, 3).Value = 0.2360468285350521
$ws.Cells.Item(24, 4).Value = 0.1053452606290648
$ws.Cells.Item(24, 5).Value = 0.1316941089510877
$ws.Cells.Item(24, 7).Value = 0.002397744461320226
$ws.Cells.Item(24, 13).Value = 0.4631092311817753
$ws.Cells.Item(24, 14).Value = 1.029113207797437
$ws.Cells.Item(24, 15).Value = 2.716405190614864
$ws.Cells.Item(25, 2).Value = 0.9271108228649041
$ws.Cells.Item(25, 3).Value = 0.1952038511880403
$ws.Cells.Item(25, 4).Value = 0.09013586295044718
$ws.Cells.Item(25, 5).Value = 0.1186029645540359
$ws.Cells.Item(25, 7).Value = 0.002403826375889016
$ws.Cells.Item(25, 13).Value = 0.3918832632879088
$ws.Cells.Item(25, 14).Value = 1.056537043870337
$ws.Cells.Item(25, 15).Value = 2.628955518964943
